$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9; everything from the old row 9 downward
# shifts down by one (old row 9 -> new row 10, ..., old row 64 -> new row 65).
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new record's data.
$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(9, 3).Value = "Ñuble"
$ws.Cells.Item(9, 4).Value = 44537
$ws.Cells.Item(9, 5).Value = 16
$ws.Cells.Item(9, 6).Value = "Fruta"
$ws.Cells.Item(9, 7).Value = 100103
$ws.Cells.Item(9, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(9, 9).Value = 100103001
$ws.Cells.Item(9, 10).Value = "Cereza"
$ws.Cells.Item(9, 11).Value = "Lapins"
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 200
$ws.Cells.Item(9, 14).Value = 12000
$ws.Cells.Item(9, 15).Value = 13000
$ws.Cells.Item(9, 16).Value = 12500
$ws.Cells.Item(9, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(9, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(9, 19).Value = 1250
$ws.Cells.Item(9, 20).Value = 10
